# Locations and prize blox flags
# Add the new Agrabah location rows (areaId/worldId/name/display) to the
# "Locations" sheet, and dummy out a handful of extra template rows for
# future entries (matching the existing pattern used throughout the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Locations")

# New Agrabah locations: B=areaId, C=worldId, D=name, E=display
$newLocations = @(
    @{ Row = 49; B = '0x6';  C = '0x4'; D = 'AgrabahPlaza';               E = 'Plaza' },
    @{ Row = 50; B = '0x8';  C = '0x4'; D = 'AgrabahAlley';               E = 'Alley' },
    @{ Row = 51; B = '0x2';  C = '0x4'; D = 'AgrabahStreet';              E = 'Main Street' },
    @{ Row = 52; B = '0x3';  C = '0x4'; D = 'AgrabahPalace';              E = 'Palace Gates' },
    @{ Row = 53; B = '0xa';  C = '0x4'; D = 'AgrabahPlazaFrozen';         E = 'Plaza' },
    @{ Row = 54; B = '0xc';  C = '0x4'; D = 'AgrabahPalaceFrozen';        E = 'Palace Gates' },
    @{ Row = 55; B = '0xb';  C = '0x4'; D = 'AgrabahStreetFrozen';        E = 'Main Street' },
    @{ Row = 56; B = '0xd';  C = '0x4'; D = 'AgrabahAlleyFrozen';         E = 'Alley' },
    @{ Row = 57; B = '0x12'; C = '0x4'; D = 'AgrabahCaveEntranceNoEntry'; E = 'Cave / Entrance' },
    @{ Row = 58; B = '0xe';  C = '0x4'; D = 'AgrabahBazaarFrozen';        E = 'Bazaar' },
    @{ Row = 59; B = '0x9';  C = '0x4'; D = 'AgrabahBazaar';              E = 'Bazaar' },
    @{ Row = 60; B = '0x1';  C = '0x4'; D = 'AgrabahCaveEntrance';        E = 'Cave / Entrance' },
    @{ Row = 61; B = '0x4';  C = '0x4'; D = 'AgrabahCaveBeginnings';      E = 'Cave / Hall of Beginnings' },
    @{ Row = 62; B = '0x5';  C = '0x4'; D = 'AgrabahGauntlet';            E = 'Cave / Gauntlet' },
    @{ Row = 63; B = '0x11'; C = '0x4'; D = 'AgrabahGauntletLower';       E = 'Cave / Gauntlet' }
)

# Every data row shares the same "key"/JSON formula in column F (and an
# incrementing index in column A); reproduce that for every row being
# added/extended below, exactly like the existing rows already do.
$formulaTemplate = '=_xlfn.CONCAT( ,A{0},": { ""worldId"": ",C{0},", ""name"": """,D{0},""", ""display"": """,E{0},""", ""areaId"": ",B{0},", },")'

foreach ($loc in $newLocations) {
    $row = $loc.Row
    $ws.Range("A$row").Value2 = $row - 1
    $ws.Range("B$row").Value2 = $loc.B
    $ws.Range("C$row").Value2 = $loc.C
    $ws.Range("D$row").Value2 = $loc.D
    $ws.Range("E$row").Value2 = $loc.E
    $ws.Range("F$row").Formula = $formulaTemplate.Replace('{0}', $row)
}

# Extend the table with extra dummied-out rows (A = index, F = the same
# CONCAT formula as every other row) ready for future entries, same as the
# two blank template rows that used to sit at the bottom of the sheet.
for ($row = 64; $row -le 72; $row++) {
    $ws.Range("A$row").Value2 = $row - 1
    $ws.Range("F$row").Formula = $formulaTemplate.Replace('{0}', $row)
}

# Move the sheet's active selection down to where the new rows were added
$ws.Activate()
[void]$ws.Range("B64").Select()
